# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" values in E16/E17 swap (2003 <-> 2002) and the
# "Salario Basico" values in G16/G17 are updated from 828116 to 877803.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" period labels between the two data rows.
$ws.Range("E16").Value = "2002"
$ws.Range("E17").Value = "2003"

# Update the "Salario Basico" amounts for both rows.
$ws.Range("G16").Value = 877803
$ws.Range("G17").Value = 877803
